$d = $word.ActiveDocument

# Update the date line (unique text in the document)
$d.Content.Find.Execute("2024-06-12 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-13 Thursday", 2) | Out-Null

$t = $d.Tables.Item(1)

# Update table cells directly by assigning Range.Text to avoid Find matching
# duplicate values elsewhere in the document (some new values equal other
# cells old values, e.g. 96÷7= and 45÷4=).
$t.Cell(1,1).Range.Text = "76÷4="
$t.Cell(1,2).Range.Text = "60÷4="
$t.Cell(1,3).Range.Text = "96÷7="
$t.Cell(1,4).Range.Text = "61÷5="
$t.Cell(1,5).Range.Text = "95÷8="

$t.Cell(5,1).Range.Text = "30÷7="
$t.Cell(5,2).Range.Text = "78÷3="
$t.Cell(5,3).Range.Text = "89÷5="
$t.Cell(5,4).Range.Text = "36÷7="
$t.Cell(5,5).Range.Text = "61÷9="

$t.Cell(9,1).Range.Text = "30÷6="
$t.Cell(9,2).Range.Text = "82÷8="
$t.Cell(9,3).Range.Text = "43÷9="
$t.Cell(9,4).Range.Text = "77÷5="
$t.Cell(9,5).Range.Text = "62÷7="

$t.Cell(13,1).Range.Text = "45÷4="
$t.Cell(13,2).Range.Text = "88÷6="
$t.Cell(13,3).Range.Text = "70÷5="
$t.Cell(13,4).Range.Text = "86÷6="
$t.Cell(13,5).Range.Text = "78÷8="

$t.Cell(17,1).Range.Text = "34÷6="
$t.Cell(17,2).Range.Text = "63÷7="
$t.Cell(17,3).Range.Text = "45÷9="
$t.Cell(17,4).Range.Text = "55÷3="
$t.Cell(17,5).Range.Text = "34÷4="
